# Updates market-board price / Leve-profit columns (H:N) for a set of
# Leve rows across the crafting-job sheets, reflecting refreshed
# currentAveragePrice data pulled by the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

# Row 100: Asking for a Friend
$ws.Range("H100").Value = 14287392
$ws.Range("I100").Value = 861.8570999999999
$ws.Range("J100").Value = 28573922
$ws.Range("K100").Value = 861.8570999999999
$ws.Range("L100").Value = 28573922
$ws.Range("M100").Value = -320.8570999999999
$ws.Range("N100").Value = -28575004

# Row 106: Making Your Mark
$ws.Range("H106").Value = 166670340
$ws.Range("I106").Value = 200003800
$ws.Range("J106").Value = 3000
$ws.Range("K106").Value = 200003800
$ws.Range("L106").Value = 3000
$ws.Range("M106").Value = -200003169
$ws.Range("N106").Value = -4262

# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 2713074.2
$ws.Range("I132").Value = 3537.7742
$ws.Range("J132").Value = 11112637
$ws.Range("K132").Value = 10613.3226
$ws.Range("L132").Value = 33337911
$ws.Range("M132").Value = -8083.3226
$ws.Range("N132").Value = -33342971


$ws = $wb.Worksheets.Item("ARM")

# Row 37: Get Shirty
$ws.Range("H37").Value = 11273.3125
$ws.Range("J37").Value = 14143.9
$ws.Range("L37").Value = 14143.9
$ws.Range("N37").Value = -14689.9

# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 4188367
$ws.Range("I61").Value = 1985346
$ws.Range("J61").Value = 19609514
$ws.Range("K61").Value = 1985346
$ws.Range("L61").Value = 19609514
$ws.Range("M61").Value = -1985134
$ws.Range("N61").Value = -19609938

# Row 80: A Squire to Inspire
$ws.Range("H80").Value = 23420.562
$ws.Range("J80").Value = 23420.562
$ws.Range("L80").Value = 23420.562
$ws.Range("N80").Value = -25416.562

# Row 83: All's Fair in Highborn Assassination (L)
$ws.Range("H83").Value = 23420.562
$ws.Range("J83").Value = 23420.562
$ws.Range("L83").Value = 70261.686
$ws.Range("N83").Value = -80245.686

# Row 136: Metal with Mettle
$ws.Range("H136").Value = 4188367
$ws.Range("I136").Value = 1985346
$ws.Range("J136").Value = 19609514
$ws.Range("K136").Value = 5956038
$ws.Range("L136").Value = 58828542
$ws.Range("M136").Value = -5953488
$ws.Range("N136").Value = -58833642


$ws = $wb.Worksheets.Item("BSM")

# Row 82: Spirituality Inspector
$ws.Range("H82").Value = 20462.285
$ws.Range("J82").Value = 27713.2
$ws.Range("L82").Value = 27713.2
$ws.Range("N82").Value = -28479.2

# Row 85: The Clamor for Hammers (L)
$ws.Range("H85").Value = 20462.285
$ws.Range("J85").Value = 27713.2
$ws.Range("L85").Value = 27713.2
$ws.Range("N85").Value = -30365.2

# Row 86: Through Thick and Thin
$ws.Range("H86").Value = 1906.98
$ws.Range("I86").Value = 1916.1414
$ws.Range("J86").Value = 1000
$ws.Range("K86").Value = 1916.1414
$ws.Range("L86").Value = 1000
$ws.Range("M86").Value = -793.1414
$ws.Range("N86").Value = -3246

# Row 89: Piercing Eyes Deserve Piercing Shafts (L)
$ws.Range("H89").Value = 1906.98
$ws.Range("I89").Value = 1916.1414
$ws.Range("J89").Value = 1000
$ws.Range("K89").Value = 9580.707
$ws.Range("L89").Value = 5000
$ws.Range("M89").Value = -3964.707
$ws.Range("N89").Value = -16232


$ws = $wb.Worksheets.Item("CRP")

# Row 31: Wall Not Found
$ws.Range("H31").Value = 2562379
$ws.Range("I31").Value = 3335297
$ws.Range("J31").Value = 1958536.9
$ws.Range("K31").Value = 3335297
$ws.Range("L31").Value = 1958536.9
$ws.Range("M31").Value = -3335002
$ws.Range("N31").Value = -1959126.9

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 2562379
$ws.Range("I34").Value = 3335297
$ws.Range("J34").Value = 1958536.9
$ws.Range("K34").Value = 3335297
$ws.Range("L34").Value = 1958536.9
$ws.Range("M34").Value = -3335095
$ws.Range("N34").Value = -1958940.9

# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 1753042
$ws.Range("I58").Value = 6364.6113
$ws.Range("J58").Value = 5683066
$ws.Range("K58").Value = 6364.6113
$ws.Range("L58").Value = 5683066
$ws.Range("M58").Value = -6161.6113
$ws.Range("N58").Value = -5683472

# Row 59: Bow Down to Magic
$ws.Range("H59").Value = 16667
$ws.Range("J59").Value = 16667
$ws.Range("L59").Value = 16667
$ws.Range("N59").Value = -18957

# Row 68: Do You Even String Bow
$ws.Range("H68").Value = 18499.727
$ws.Range("J68").Value = 17687.25
$ws.Range("L68").Value = 17687.25
$ws.Range("N68").Value = -19185.25

# Row 71: Win One Bow, Get Three Free (L)
$ws.Range("H71").Value = 18499.727
$ws.Range("J71").Value = 17687.25
$ws.Range("L71").Value = 53061.75
$ws.Range("N71").Value = -60549.75

# Row 80: The Long Armillae of the Law
$ws.Range("H80").Value = 21300
$ws.Range("J80").Value = 21300
$ws.Range("L80").Value = 21300
$ws.Range("N80").Value = -23546

# Row 83: Wooden Ambitions (L)
$ws.Range("H83").Value = 21300
$ws.Range("J83").Value = 21300
$ws.Range("L83").Value = 63900
$ws.Range("N83").Value = -75132

# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 2679.7896
$ws.Range("I132").Value = 2075.5
$ws.Range("J132").Value = 3119.2727
$ws.Range("K132").Value = 6226.5
$ws.Range("L132").Value = 9357.8181
$ws.Range("M132").Value = -3696.5
$ws.Range("N132").Value = -14417.8181

# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 1826088.9
$ws.Range("I134").Value = 10018.538
$ws.Range("K134").Value = 30055.614
$ws.Range("M134").Value = -27520.614

# Row 136: Turali Quality
$ws.Range("H136").Value = 1753042
$ws.Range("I136").Value = 6364.6113
$ws.Range("J136").Value = 5683066
$ws.Range("K136").Value = 19093.8339
$ws.Range("L136").Value = 17049198
$ws.Range("M136").Value = -16543.8339
$ws.Range("N136").Value = -17054298


$ws = $wb.Worksheets.Item("CUL")

# Row 5: What a Sap
$ws.Range("H5").Value = 3984362.5
$ws.Range("J5").Value = 3907385
$ws.Range("L5").Value = 11722155
$ws.Range("N5").Value = -11722379

# Row 135: Not-so-secret Ingredient
$ws.Range("H135").Value = 3984362.5
$ws.Range("J135").Value = 3907385
$ws.Range("L135").Value = 35166465
$ws.Range("N135").Value = -35171535


$ws = $wb.Worksheets.Item("GSM")

# Row 70: Sky Is the Limit
$ws.Range("H70").Value = 7635481.5
$ws.Range("I70").Value = 3972670.8
$ws.Range("J70").Value = 11908761
$ws.Range("K70").Value = 3972670.8
$ws.Range("L70").Value = 11908761
$ws.Range("M70").Value = -3972400.8
$ws.Range("N70").Value = -11909301

# Row 73: Hulls of Broken Dreams (L)
$ws.Range("H73").Value = 7635481.5
$ws.Range("I73").Value = 3972670.8
$ws.Range("J73").Value = 11908761
$ws.Range("K73").Value = 3972670.8
$ws.Range("L73").Value = 11908761
$ws.Range("M73").Value = -3971734.8
$ws.Range("N73").Value = -11910633

# Row 132: On Board for Lar
$ws.Range("H132").Value = 8910226
$ws.Range("I132").Value = 8255248.5
$ws.Range("K132").Value = 24765745.5
$ws.Range("M132").Value = -24763215.5


$ws = $wb.Worksheets.Item("LTW")

# Row 7: Tan Before the Ban
$ws.Range("H7").Value = 1327.8889
$ws.Range("I7").Value = 1327.8889
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 1327.8889
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -1215.8889
$ws.Range("N7").Value = ""

# Row 126: Battered Books
$ws.Range("H126").Value = 1327.8889
$ws.Range("I126").Value = 1327.8889
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 3983.6667
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -1513.6667
$ws.Range("N126").Value = ""


$ws = $wb.Worksheets.Item("WVR")

# Row 4: Not Cool Enough
$ws.Range("H4").Value = 57833.332
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 57833.332
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 57833.332
$ws.Range("N4").Value = -58059.332
$ws.Range("M4").Value = ""
